$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear A3 (was "test 3") entirely - the string is no longer used anywhere
# so it drops out of the shared strings table, shifting later indices down.
$ws.Range("A3").ClearContents()

# Update selection to A3
$ws.Range("A3").Select()
